$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new client code (20129) to the group stored in B2
$ws.Range("B2").Value = "233.10167.20234.20292.20357.20379.20385.50102.50818.60124.60139.60158.60174.60258.70101.3.70114.20129"

# Move the active selection to B3 (as left by the editor after the change)
$ws.Range("B3").Select()
